$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28, shifting existing rows 28-48 down to 29-49
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with the new weekly price record
$ws.Range("A28").Value = 3
$ws.Range("B28").Value = "Femacal de La Calera"
$ws.Range("C28").Value = "Coquimbo"
$ws.Range("D28").Value = 44588
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100108
$ws.Range("H28").Value = "Tropicales y subtropicales"
$ws.Range("I28").Value = 100108004
$ws.Range("J28").Value = "Papaya"
$ws.Range("K28").Value = "Cultivar IV Región"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 50
$ws.Range("N28").Value = 34000
$ws.Range("O28").Value = 34000
$ws.Range("P28").Value = 34000
$ws.Range("Q28").Value = "`$/caja 15 kilos granel"
$ws.Range("R28").Value = "Provincia del Elquí"
$ws.Range("S28").Value = 2267
$ws.Range("T28").Value = 15
